# Add new kenmerktypes (codelist entries) for trajecten: actual/optimal
# time and speed, to the "VkmVerkeersKenmerkType" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # VkmVerkeersKenmerkType

# Columns: A=Klasse, B=Definitie, C=Notation, D=Label, E=Status

$data = @(
    @('actuele_tijd', 'de huidige tijd gemeten om een traject af te leggen', 'actualTt', 'actuele tijd', 'ingebruik'),
    @('vertraging', 'de tijd waarin de actuele tijd afwijkt van de optimale tijd om een traject af te leggen', 'delay', 'vertraging', 'ingebruik'),
    @('optimale_tijd', 'de optimale tijd om een traject af te leggen', 'optimalTt', 'optimale tijd', 'ingebruik'),
    @('trajectgemiddelde snelheid', 'gemiddelde snelheid gemeten over een traject (harmonisch gemiddelde)', 'actualKmH', 'trajectgemiddelde snelheid', 'ingebruik'),
    @('optimale_trajectgemiddelde_snelheid', 'de optimale gemiddelde snelheid om een traject af te leggen', 'optimalKmH', 'optimale trajectgemiddelde snelheid', 'ingebruik')
)

$startRow = 8
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Rows.Item($r).RowHeight = 17.25

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Enable iterative calculation (matches calcPr iterateDelta="1E-4")
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# The edit session ends with this sheet active/selected, scrolled so the
# next empty row is in view.
$ws.Activate()
$ws.Range("B13").Select()
